$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (D1:F1) to match the new sharedStrings
# entries (ORG_SEC_IDENOLD, ORG_SEC_IDENNEW, ORG_SEC_STATUS).
$ws.Range("D1").Value = "ORG_SEC_IDENOLD"
$ws.Range("E1").Value = "ORG_SEC_IDENNEW"
$ws.Range("F1").Value = "ORG_SEC_STATUS"

# Move the active selection, matching the post-edit workbook state.
$ws.Range("E7").Select()
